$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 5: Runmode changed from N to Y
$ws.Range("D5").Value = "Y"

# Row 7: new test case details
$ws.Range("C7").Value = 'Verify that the "Thanks for your interest in EndNote......" modal displayed when user clicks on the export button when user is signed to facebook account and not having existing steam account'
$ws.Range("B7").Value = "OPQA-1701"
$ws.Range("D7").Value = "Y"
$ws.Rows.Item(7).RowHeight = 45

# Update selection/view
$ws.Range("D4").Select()
